# Edit script: restructure progress-report body paragraphs and
# remove the unwanted proofErr split around "Thoa" in the header.

$d = $word.ActiveDocument

# --- 1. Rebuild the body content -------------------------------------------------
# Delete all existing body content (this also clears paragraph formatting,
# e.g. the firstLine indent, matching the target which has no <w:pPr>).
$d.Content.Delete()

$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>First we filled in the Node class. Each node has a data field and a link to the next Node (one way link).</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">Then we filled out the code for </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>ListQueue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">. </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">We made the constructors, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>enqueue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>dequeue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">, and the front methods. For </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">each we figured out the types of cases there would be like an empty queue, a queue with one element, and the regular case. </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">We have </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>testEmpty</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve">) and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>testOne</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">() to work. And we are now working on the remaining test: </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>testMany</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve">). </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">For each test, we test the </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>front(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve">), then </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>dequeue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">(), then </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>enqueue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">(), in the order of how those methods impact the queue’s original state. For example, when we want to test the methods for an empty queue, if we </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>enqueue</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t>) before front(), then the queue is no longer empty.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$d.Content.InsertXML($bodyXml)

# --- 2. Fix up the header: merge " & " + "Thoa" + " Ta" into a single run --------
# (visible text is unchanged, only the spell-check proofErr split around "Thoa"
# is removed, i.e. the three runs become one plain run)
$hdr = $d.Sections.Item(1).Headers.Item(1)
$find = $hdr.Range.Find
$find.ClearFormatting()
$found = $find.Execute(" & Thoa Ta", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find header text ' & Thoa Ta' to fix up"
}
$matched = $find.Parent
$matched.Text = ""
$matched.InsertAfter(" & Thoa Ta")

Write-Host "Body paragraphs:" $d.Paragraphs.Count
Write-Host "Body text:" $d.Content.Text
Write-Host "Header text:" $d.Sections.Item(1).Headers.Item(1).Range.Text
